{"js": "// Replace each two-digit multiplication expression in the table with its\n// updated counterpart, per the commit's data update. Each \"NNxNN=\" string\n// is unique in the document, so an exact, case-sensitive whole-document\n// search safely targets the single matching run.\nconst replacements = [\n  [\"45\u00d775=\", \"44\u00d795=\"],\n  [\"31\u00d774=\", \"75\u00d733=\"],\n  [\"13\u00d766=\", \"78\u00d781=\"],\n  [\"92\u00d742=\", \"39\u00d729=\"],\n  [\"22\u00d714=\", \"27\u00d714=\"],\n  [\"91\u00d764=\", \"18\u00d728=\"],\n  [\"45\u00d785=\", \"62\u00d795=\"],\n  [\"84\u00d743=\", \"46\u00d735=\"],\n  [\"14\u00d756=\", \"65\u00d715=\"],\n  [\"47\u00d741=\", \"15\u00d736=\"],\n  [\"37\u00d793=\", \"32\u00d717=\"],\n  [\"83\u00d758=\", \"52\u00d778=\"],\n  [\"76\u00d769=\", \"37\u00d766=\"],\n  [\"61\u00d755=\", \"11\u00d776=\"],\n  [\"79\u00d749=\", \"70\u00d750=\"],\n  [\"97\u00d753=\", \"82\u00d725=\"],\n  [\"48\u00d723=\", \"74\u00d734=\"],\n  [\"37\u00d773=\", \"15\u00d720=\"],\n  [\"55\u00d788=\", \"58\u00d798=\"],\n  [\"42\u00d749=\", \"81\u00d752=\"],\n  [\"17\u00d775=\", \"84\u00d771=\"],\n  [\"53\u00d740=\", \"77\u00d762=\"],\n  [\"31\u00d725=\", \"76\u00d761=\"],\n  [\"23\u00d792=\", \"79\u00d719=\"],\n  [\"78\u00d713=\", \"73\u00d784=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n\n", "ps1": "# Update each two-digit multiplication expression in the practice-sheet\n# table to its new value (per the data refresh in this commit). Each\n# \"NNxNN=\" string occurs exactly once in the document, so a plain\n# Find & Replace targeting the whole document body is safe for every\n# pair and will not clobber an already-updated cell.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"45\u00d775=\", \"44\u00d795=\"),\n  @(\"31\u00d774=\", \"75\u00d733=\"),\n  @(\"13\u00d766=\", \"78\u00d781=\"),\n  @(\"92\u00d742=\", \"39\u00d729=\"),\n  @(\"22\u00d714=\", \"27\u00d714=\"),\n  @(\"91\u00d764=\", \"18\u00d728=\"),\n  @(\"45\u00d785=\", \"62\u00d795=\"),\n  @(\"84\u00d743=\", \"46\u00d735=\"),\n  @(\"14\u00d756=\", \"65\u00d715=\"),\n  @(\"47\u00d741=\", \"15\u00d736=\"),\n  @(\"37\u00d793=\", \"32\u00d717=\"),\n  @(\"83\u00d758=\", \"52\u00d778=\"),\n  @(\"76\u00d769=\", \"37\u00d766=\"),\n  @(\"61\u00d755=\", \"11\u00d776=\"),\n  @(\"79\u00d749=\", \"70\u00d750=\"),\n  @(\"97\u00d753=\", \"82\u00d725=\"),\n  @(\"48\u00d723=\", \"74\u00d734=\"),\n  @(\"37\u00d773=\", \"15\u00d720=\"),\n  @(\"55\u00d788=\", \"58\u00d798=\"),\n  @(\"42\u00d749=\", \"81\u00d752=\"),\n  @(\"17\u00d775=\", \"84\u00d771=\"),\n  @(\"53\u00d740=\", \"77\u00d762=\"),\n  @(\"31\u00d725=\", \"76\u00d761=\"),\n  @(\"23\u00d792=\", \"79\u00d719=\"),\n  @(\"78\u00d713=\", \"73\u00d784=\"),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        throw \"Could not find text to replace: $oldText\"\n    }\n}\n\n"}
